# Applies the "Updated cryptos list" data refresh (Tue May 14 22:39:49 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.610.85"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "2.890.70"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.63"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.62"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.500"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Value = "2.887.50"
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.58"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "3.370.24"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").Value = "61.580.38"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "2.892.21"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.10"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.98"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.650"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.79"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.85"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.91"
$ws.Range("E27").Value = "  -11.34%  "
$ws.Range("E28").Value = "  -5.67%  "
$ws.Range("E29").Value = "  +7.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.98"
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("E32").Value = "  -8.92%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.40"
$ws.Range("E35").Value = "  -3.51%  "
$ws.Range("E36").Value = "  -3.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.34"
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.78"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.79"
$ws.Range("E40").Value = "  -7.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.14"
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("E42").Value = "  -4.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.25"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.266"
$ws.Range("E44").Value = "  -4.49%  "
$ws.Range("D45").Value = "2.678.26"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.02"
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "342.44"
$ws.Range("E48").Value = "  -3.81%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.44"
$ws.Range("E51").Value = "  -5.43%  "
